# Update "想去人数" (column F) values across the four worksheets to match
# the newly-scraped counts.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 392
$ws1.Range("F5").Value  = 910
$ws1.Range("F6").Value  = 145
$ws1.Range("F7").Value  = 915
$ws1.Range("F9").Value  = 166
$ws1.Range("F12").Value = 758
$ws1.Range("F14").Value = 539
$ws1.Range("F16").Value = 1274
$ws1.Range("F19").Value = 1030
$ws1.Range("F20").Value = 2766
$ws1.Range("F21").Value = 1240
$ws1.Range("F22").Value = 631
$ws1.Range("F24").Value = 1229
$ws1.Range("F25").Value = 51
$ws1.Range("F26").Value = 956
$ws1.Range("F27").Value = 311
$ws1.Range("F28").Value = 340
$ws1.Range("F29").Value = 1291

# 演出 (Performance) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 509
$ws2.Range("F4").Value = 346

# 本地生活 (Local life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 714

# 全部类型 (All types) sheet - aggregate of all three above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 714
$ws4.Range("F5").Value  = 392
$ws4.Range("F7").Value  = 509
$ws4.Range("F8").Value  = 509
$ws4.Range("F9").Value  = 346
$ws4.Range("F12").Value = 910
$ws4.Range("F13").Value = 145
$ws4.Range("F15").Value = 915
$ws4.Range("F17").Value = 166
$ws4.Range("F25").Value = 758
$ws4.Range("F27").Value = 539
$ws4.Range("F29").Value = 1274
$ws4.Range("F32").Value = 1030
$ws4.Range("F33").Value = 2766
$ws4.Range("F34").Value = 1240
$ws4.Range("F35").Value = 631
$ws4.Range("F37").Value = 1229
$ws4.Range("F38").Value = 51
$ws4.Range("F41").Value = 311
$ws4.Range("F42").Value = 340
$ws4.Range("F43").Value = 1291

$wb.Save()
